$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "ILO" source note text before we overwrite A20/A21.
$iloName = $ws.Range("A20").Text
$iloText = $ws.Range("A21").Text

# ---------------------------------------------------------------
# New "Number of employees / Assets / Turnover" breakdown table
# inserted at rows 17-21 (pushing the ILO source note down to 26-27)
# ---------------------------------------------------------------
$ws.Range("B17").Value = "Number of employees"
$ws.Range("C17").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D17").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B17:D17").Style = "title"
$ws.Range("B17:D17").Font.Bold = $true

$ws.Range("A18").Value = "Micro"
$ws.Range("B18").Value = "<5"
$ws.Range("C18").Value = "'"
$ws.Range("D18").Value = "'"
$ws.Range("A18:D18").Style = "Normal"

$ws.Range("A19").Value = "Small"
$ws.Range("B19").Value = "5-9"
$ws.Range("C19").Value = "'"
$ws.Range("D19").Value = "'"
$ws.Range("A19:D19").Style = "Normal"

$ws.Range("A20").Value = "Medium"
$ws.Range("B20").Value = "10-29"
$ws.Range("C20").Value = "'"
$ws.Range("D20").Value = "'"
$ws.Range("A20:D20").Style = "Normal"

$ws.Range("A21").Value = "Large"
$ws.Range("B21").Value = ">=30"
$ws.Range("C21").Value = "'"
$ws.Range("D21").Value = "'"
$ws.Range("A21:D21").Style = "Normal"

# ---------------------------------------------------------------
# Move the existing ILO source note down to rows 26-27
# ---------------------------------------------------------------
$ws.Range("A26").Value = $iloName
$ws.Range("A26").Style = "title"
$ws.Range("A26").Font.Bold = $true

$ws.Range("A27").Value = $iloText
$ws.Range("A27").Style = "source"
$ws.Range("A27").Font.Italic = $true
